# EPBDS-9540 Support Java Name convention on Json field name generating in SpreadsheetResults. Rework.
# The step names used inside the _res_.$Step2["..."] / ..."]:Integer text labels
# are normalized to match the Java-bean-style capitalization of the actual
# step names (Step1 / Step2 / Step3 / SomeStep) instead of the lower-camel
# variants that were there before (step1 / step2 / step3 / someStep).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D29").Value = '_res_.$Step2["Step1"]:Integer'
$ws.Range("E29").Value = '_res_.$Step2["SomeStep"]:Integer'
$ws.Range("F29").Value = '_res_.$Step2["Step2"]:Integer'

$ws.Range("D30").Value = '_res_.$Step2["Step1"]'
$ws.Range("E30").Value = '_res_.$Step2["SomeStep"]'
$ws.Range("F30").Value = '_res_.$Step2["Step2"]'

$ws.Range("C45").Value = '_res_.$Step2["Step1"]:Integer'
$ws.Range("D45").Value = '_res_.$Step2["Step2"]:Integer'
$ws.Range("E45").Value = '_res_.$Step2["Step3"]:Integer'

$ws.Range("C46").Value = '_res_.$Step2["Step1"]:Integer'
$ws.Range("D46").Value = '_res_.$Step2["Step2"]:Integer'
$ws.Range("E46").Value = '_res_.$Step2["Step3"]:Integer'
